# Sara-Alert-Format-Invalid-Monitorees.xlsx:
# "reorder race options for consistency" -- the "Race Unknown" and
# "Race Other" column headers (CV1/CW1 on Sheet1) swap places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$raceUnknownCol = $ws.Range("CV1").Value2
$raceOtherCol   = $ws.Range("CW1").Value2

$ws.Range("CV1").Value = $raceOtherCol
$ws.Range("CW1").Value = $raceUnknownCol
